$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Remove the now-unused header_4 / row_1_d / row_2_d / sum_4 / row_3_c / row_3_d
# columns (C and D) so the shared-string table gets compacted on save.
$ws.Range("C1:D4").ClearContents()
$ws.Range("C7:D7").ClearContents()

# Push the summary row (currently on row 7) down to row 12, opening up
# rows 5-9 for the new data rows (row_4 .. row_8).
$ws.Range("A5:A9").EntireRow.Insert()

# Fill in the new rows. The order below matches how the shared strings
# table ends up populated (note row_4_b / B5 is filled in last).
$ws.Range("A5").Value = "row_4_a"
$ws.Range("A6").Value = "row_5_a"
$ws.Range("B6").Value = "row_5_b"
$ws.Range("A7").Value = "row_6_a"
$ws.Range("B7").Value = "row_6_b"
$ws.Range("A8").Value = "row_7_a"
$ws.Range("B8").Value = "row_7_b"
$ws.Range("A9").Value = "row_8_a"
$ws.Range("B9").Value = "row_8_b"
$ws.Range("B5").Value = "row_4_b"

# Update the selection to match the final state of the worksheet.
$ws.Range("B6").Select()
